$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 55 (shifts existing rows 55-73 down to 56-74)
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new record
$ws.Range("A55").Value = 4
$ws.Range("B55").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C55").Value = "Los Lagos"
$ws.Range("D55").Value = 45215
$ws.Range("E55").Value = 10
$ws.Range("F55").Value = 300000000
$ws.Range("G55").Value = "Espárragos"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 2000
$ws.Range("M55").Value = 2000
$ws.Range("N55").Value = "$/kilo"
$ws.Range("O55").Value = "Provincia de Linares"
$ws.Range("P55").Value = 2000
$ws.Range("Q55").Value = 1
$ws.Range("R55").Value = "Hortaliza"
